$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.03158659211277
$ws.Range("D2").Value = 1.039944605764331
$ws.Range("E2").Value = 1.049178480035779
$ws.Range("F2").Value = 1.05359276745835
$ws.Range("I2").Value = 1.034747391137029
$ws.Range("J2").Value = 1.03672150028258
$ws.Range("K2").Value = 1.042728354770345
$ws.Range("L2").Value = 1.051936261370192
$ws.Range("M2").Value = 1.056338309151664
$ws.Range("N2").Value = 1.0381937637266

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.03254069711127
$ws.Range("D3").Value = 1.040678811066486
$ws.Range("E3").Value = 1.050126938082122
$ws.Range("F3").Value = 1.054531783510816
$ws.Range("I3").Value = 1.034919662780131
$ws.Range("J3").Value = 1.037317519060715
$ws.Range("K3").Value = 1.043273082932233
$ws.Range("L3").Value = 1.052696561341941
$ws.Range("M3").Value = 1.057090078103249
$ws.Range("N3").Value = 1.038790628919763

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.033158193699753
$ws.Range("D4").Value = 1.041153666245438
$ws.Range("E4").Value = 1.050741608727319
$ws.Range("F4").Value = 1.0551399192101
$ws.Range("I4").Value = 1.035029272054694
$ws.Range("J4").Value = 1.03770270644114
$ws.Range("K4").Value = 1.043624667336109
$ws.Range("L4").Value = 1.053188856267989
$ws.Range("M4").Value = 1.057576425525737
$ws.Range("N4").Value = 1.039176363310451

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.033417818847297
$ws.Range("D5").Value = 1.041353240220676
$ws.Range("E5").Value = 1.051000243613344
$ws.Range("F5").Value = 1.055395704937971
$ws.Range("I5").Value = 1.035074905334807
$ws.Range("J5").Value = 1.037864524234619
$ws.Range("K5").Value = 1.04377225904767
$ws.Range("L5").Value = 1.053395894987299
$ws.Range("M5").Value = 1.05778086154156
$ws.Range("N5").Value = 1.039338410903757

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.033461412763755
$ws.Range("D6").Value = 1.041386746296722
$ws.Range("E6").Value = 1.051043682844095
$ws.Range("F6").Value = 1.055438659807854
$ws.Range("I6").Value = 1.035082541162531
$ws.Range("J6").Value = 1.037891687386065
$ws.Range("K6").Value = 1.043797027744533
$ws.Range("L6").Value = 1.053430662238091
$ws.Range("M6").Value = 1.057815185781967
$ws.Range("N6").Value = 1.039365612629993

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.033161662708435
$ws.Range("D7").Value = 1.04115633318103
$ws.Range("E7").Value = 1.050745063727582
$ws.Range("F7").Value = 1.055143336540092
$ws.Range("I7").Value = 1.035029883563975
$ws.Range("J7").Value = 1.037704869111325
$ws.Range("K7").Value = 1.043626640307426
$ws.Range("L7").Value = 1.05319162242462
$ws.Range("M7").Value = 1.057579157307253
$ws.Range("N7").Value = 1.039178529051875

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.031909009624928
$ws.Range("D8").Value = 1.040192779617166
$ws.Range("E8").Value = 1.04949881741044
$ws.Range("F8").Value = 1.053910001914703
$ws.Range("I8").Value = 1.03480599626939
$ws.Range("J8").Value = 1.036923025781549
$ws.Range("K8").Value = 1.042912632102027
$ws.Range("L8").Value = 1.052193139332416
$ws.Range("M8").Value = 1.056592392467593
$ws.Range("N8").Value = 1.038395575414893

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.029702680796853
$ws.Range("D9").Value = 1.038493212451244
$ws.Range("E9").Value = 1.047310140310284
$ws.Range("F9").Value = 1.051740820829437
$ws.Range("I9").Value = 1.034397247841781
$ws.Range("J9").Value = 1.03554170593452
$ws.Range("K9").Value = 1.041647678442359
$ws.Range("L9").Value = 1.050436259064069
$ws.Range("M9").Value = 1.054852890782191
$ws.Range("N9").Value = 1.037012293935241

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.028232515488964
$ws.Range("D10").Value = 1.037359127439234
$ws.Range("E10").Value = 1.045856049028846
$ws.Range("F10").Value = 1.050297539638807
$ws.Range("I10").Value = 1.03411522110966
$ws.Range("J10").Value = 1.034618447900607
$ws.Range("K10").Value = 1.040799877029423
$ws.Range("L10").Value = 1.049266798793141
$ws.Range("M10").Value = 1.053692820679254
$ws.Range("N10").Value = 1.036087724769021

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.02759609723092
$ws.Range("D11").Value = 1.036867825327641
$ws.Range("E11").Value = 1.045227618703731
$ws.Range("F11").Value = 1.04967327174647
$ws.Range("I11").Value = 1.033990848692951
$ws.Range("J11").Value = 1.034218113108626
$ws.Range("K11").Value = 1.04043171468871
$ws.Range("L11").Value = 1.048760847679693
$ws.Range("M11").Value = 1.053190415493861
$ws.Range("N11").Value = 1.035686821455716

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.027359729762154
$ws.Range("D12").Value = 1.036685299655391
$ws.Range("E12").Value = 1.044994373273416
$ws.Range("F12").Value = 1.049441494559167
$ws.Range("I12").Value = 1.033944313359768
$ws.Range("J12").Value = 1.034069327879855
$ws.Range("K12").Value = 1.040294804661444
$ws.Range("L12").Value = 1.04857298093429
$ws.Range("M12").Value = 1.053003787769152
$ws.Range("N12").Value = 1.035537824934854

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.027410430172615
$ws.Range("D13").Value = 1.036724453559204
$ws.Range("E13").Value = 1.045044396974691
$ws.Range("F13").Value = 1.049491206840881
$ws.Range("I13").Value = 1.033954310626814
$ws.Range("J13").Value = 1.034101246562952
$ws.Range("K13").Value = 1.040324179468541
$ws.Range("L13").Value = 1.048613275985283
$ws.Range("M13").Value = 1.053043820579318
$ws.Range("N13").Value = 1.035569788946142

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.027576558475853
$ws.Range("D14").Value = 1.03685273838903
$ws.Range("E14").Value = 1.045208334852914
$ws.Range("F14").Value = 1.049654110849518
$ws.Range("I14").Value = 1.033987008956958
$ws.Range("J14").Value = 1.034205816165221
$ws.Range("K14").Value = 1.040420400883875
$ws.Range("L14").Value = 1.048745317203308
$ws.Range("M14").Value = 1.053174989028194
$ws.Range("N14").Value = 1.035674507049241

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.027678919150696
$ws.Range("D15").Value = 1.036931774412607
$ws.Range("E15").Value = 1.04530936648543
$ws.Range("F15").Value = 1.049754495162775
$ws.Range("I15").Value = 1.034007110722993
$ws.Range("J15").Value = 1.034270233964751
$ws.Range("K15").Value = 1.040479665145758
$ws.Range("L15").Value = 1.048826680934169
$ws.Range("M15").Value = 1.053255804670045
$ws.Range("N15").Value = 1.035739016329435

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.028274756170571
$ws.Range("D16").Value = 1.037391728679579
$ws.Range("E16").Value = 1.045897781278905
$ws.Range("F16").Value = 1.050338984720973
$ws.Range("I16").Value = 1.034123427881367
$ws.Range("J16").Value = 1.034645005111998
$ws.Range("K16").Value = 1.040824288523913
$ws.Range("L16").Value = 1.049300386307266
$ws.Range("M16").Value = 1.053726161907907
$ws.Range("N16").Value = 1.036114319694698

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.028648555748163
$ws.Range("D17").Value = 1.037680183579206
$ws.Range("E17").Value = 1.046267200644933
$ws.Range("F17").Value = 1.050705803109232
$ws.Range("I17").Value = 1.034195787803326
$ws.Range("J17").Value = 1.034879940309147
$ws.Range("K17").Value = 1.041040178863142
$ws.Range("L17").Value = 1.049597645770194
$ws.Range("M17").Value = 1.054021181964009
$ws.Range("N17").Value = 1.036349588526774

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.028866603370556
$ws.Range("D18").Value = 1.037848411606394
$ws.Range("E18").Value = 1.046482792492222
$ws.Range("F18").Value = 1.050919827830423
$ws.Range("I18").Value = 1.034237776701285
$ws.Range("J18").Value = 1.035016920147556
$ws.Range("K18").Value = 1.041166001831443
$ws.Range("L18").Value = 1.049771073888297
$ws.Range("M18").Value = 1.054193253691487
$ws.Range("N18").Value = 1.036486762892266

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.02894095478087
$ws.Range("D19").Value = 1.037905769105259
$ws.Range("E19").Value = 1.046556323343703
$ws.Range("F19").Value = 1.050992815847702
$ws.Range("I19").Value = 1.034252056933679
$ws.Range("J19").Value = 1.03506361756898
$ws.Range("K19").Value = 1.041208886846701
$ws.Range("L19").Value = 1.049830215438209
$ws.Range("M19").Value = 1.054251924234163
$ws.Range("N19").Value = 1.036533526629385

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.028608448864074
$ws.Range("D20").Value = 1.03764923745138
$ws.Range("E20").Value = 1.046227553458838
$ws.Range("F20").Value = 1.050666440157095
$ws.Range("I20").Value = 1.034188046752898
$ws.Range("J20").Value = 1.034854739554118
$ws.Range("K20").Value = 1.041017026438669
$ws.Range("L20").Value = 1.049565748322385
$ws.Range("M20").Value = 1.053989529973665
$ws.Range("N20").Value = 1.036324351983783

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.027527637091727
$ws.Range("D21").Value = 1.036814962638057
$ws.Range("E21").Value = 1.045160054212146
$ws.Range("F21").Value = 1.04960613680705
$ws.Range("I21").Value = 1.033977389436019
$ws.Range("J21").Value = 1.034175025310464
$ws.Range("K21").Value = 1.040392070433422
$ws.Range("L21").Value = 1.048706432525095
$ws.Range("M21").Value = 1.053136363514223
$ws.Range("N21").Value = 1.035643672467938

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.026848242087203
$ws.Range("D22").Value = 1.036290223728967
$ws.Range("E22").Value = 1.044489926170466
$ws.Range("F22").Value = 1.048940082586638
$ws.Range("I22").Value = 1.033842986247278
$ws.Range("J22").Value = 1.033747181793009
$ws.Range("K22").Value = 1.03999822128297
$ws.Range("L22").Value = 1.048166529068174
$ws.Range("M22").Value = 1.052599874576435
$ws.Range("N22").Value = 1.035215221363614

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.027208387407251
$ws.Range("D23").Value = 1.036568416036386
$ws.Range("E23").Value = 1.044845073702828
$ws.Range("F23").Value = 1.049293113209149
$ws.Range("I23").Value = 1.033914420990646
$ws.Range("J23").Value = 1.033974034998489
$ws.Range("K23").Value = 1.040207094547773
$ws.Range("L23").Value = 1.048452705521222
$ws.Range("M23").Value = 1.052884283741753
$ws.Range("N23").Value = 1.035442396726667

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.028626571387131
$ws.Range("D24").Value = 1.037663220745132
$ws.Range("E24").Value = 1.046245467957982
$ws.Range("F24").Value = 1.050684226377691
$ws.Range("I24").Value = 1.0341915452722
$ws.Range("J24").Value = 1.034866126857093
$ws.Range("K24").Value = 1.041027488338577
$ws.Range("L24").Value = 1.049580161277251
$ws.Range("M24").Value = 1.054003832171919
$ws.Range("N24").Value = 1.036335755458034

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.030272945909832
$ws.Range("D25").Value = 1.038932780308824
$ws.Range("E25").Value = 1.047875085123583
$ws.Range("F25").Value = 1.052301111108274
$ws.Range("I25").Value = 1.034504601350716
$ws.Range("J25").Value = 1.035899233060412
$ws.Range("K25").Value = 1.041975496593005
$ws.Range("L25").Value = 1.050890143204972
$ws.Range("M25").Value = 1.055302669362153
$ws.Range("N25").Value = 1.037370328790661

